$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "ActiveSheet:" $ws.Name
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
